# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed "K" (s_vals) values for rows 2..43 (column G)
$newK = @{
    2  = 0
    3  = 5
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 3
    11 = 2
    12 = 2
    13 = 0
    14 = 3
    15 = 0
    16 = 2
    17 = 2
    18 = 1
    19 = 0
    20 = 3
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 2
    34 = 0
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
